# Auto-generated edit script applying the diff to rows 5-34
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 79244
$ws.Range("A7").Value = 130961962
$ws.Range("B7").Value = 79244
$ws.Range("Q7").Value = 446084
$ws.Range("R7").Value = 6759981
$ws.Range("A8").Value = 130961458
$ws.Range("B8").Value = 79244
$ws.Range("Q8").Value = 446059
$ws.Range("R8").Value = 6760088
$ws.Range("B9").Value = 79244
$ws.Range("A11").Value = 130961218
$ws.Range("AC11").ClearContents()
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("M11").Value = "äldre spår"
$ws.Range("A12").Value = 130961179
$ws.Range("AC12").Value = "Ringhack på stam i bakgrund"
$ws.Range("B12").Value = 79863
$ws.Range("E12").Value = 6453
$ws.Range("F12").Value = "Vedskivlav"
$ws.Range("G12").Value = "Hertelidea botryosa"
$ws.Range("H12").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M12").ClearContents()
$ws.Range("B13").Value = 79244
$ws.Range("B14").Value = 79244
$ws.Range("B15").Value = 79244
$ws.Range("B16").Value = 79244
$ws.Range("A17").Value = 130960789
$ws.Range("B17").Value = 79244
$ws.Range("Q17").Value = 446284
$ws.Range("R17").Value = 6759886
$ws.Range("A18").Value = 130960378
$ws.Range("B18").Value = 57884
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("M18").Value = "äldre spår"
$ws.Range("Q18").Value = 446272
$ws.Range("R18").Value = 6759739
$ws.Range("B19").Value = 79863
$ws.Range("A20").Value = 130960843
$ws.Range("B20").Value = 79244
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("M20").ClearContents()
$ws.Range("Q20").Value = 446247
$ws.Range("R20").Value = 6759903
$ws.Range("B21").Value = 79244
$ws.Range("A22").Value = 130962722
$ws.Range("AB22").Value = "10:26"
$ws.Range("AC22").ClearContents()
$ws.Range("B22").Value = 79863
$ws.Range("E22").Value = 6453
$ws.Range("F22").Value = "Vedskivlav"
$ws.Range("G22").Value = "Hertelidea botryosa"
$ws.Range("H22").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q22").Value = 446008
$ws.Range("R22").Value = 6759948
$ws.Range("Z22").Value = "10:26"
$ws.Range("A23").Value = 130962640
$ws.Range("B23").Value = 79863
$ws.Range("Q23").Value = 446038
$ws.Range("R23").Value = 6759945
$ws.Range("A24").Value = 130963976
$ws.Range("AB24").Value = "14:08"
$ws.Range("AC24").Value = "Miljöbild"
$ws.Range("B24").Value = 79244
$ws.Range("E24").Value = 6425
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("Q24").Value = 445929
$ws.Range("R24").Value = 6760099
$ws.Range("Z24").Value = "14:08"
$ws.Range("B25").Value = 79244
$ws.Range("A27").Value = 130961461
$ws.Range("B27").Value = 79244
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 446088
$ws.Range("R27").Value = 6760088
$ws.Range("B28").Value = 79244
$ws.Range("A29").Value = 130963807
$ws.Range("AB29").Value = "14:08"
$ws.Range("B29").Value = 57881
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 445932
$ws.Range("R29").Value = 6760079
$ws.Range("Z29").Value = "14:08"
$ws.Range("A30").Value = 130962736
$ws.Range("AB30").Value = "10:26"
$ws.Range("B30").Value = 79834
$ws.Range("E30").Value = 229821
$ws.Range("F30").Value = "Vedflamlav"
$ws.Range("G30").Value = "Ramboldia elabens"
$ws.Range("H30").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M30").ClearContents()
$ws.Range("Q30").Value = 446008
$ws.Range("R30").Value = 6759948
$ws.Range("Z30").Value = "10:26"
$ws.Range("B31").Value = 79244
$ws.Range("B32").Value = 79244
$ws.Range("B33").Value = 79244
$ws.Range("B34").Value = 79244
